$d = $word.ActiveDocument

# 1) Remove the _GoBack bookmark that was on the "Comparison of transactions scanned" paragraph
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Insert the new "Time comparison" section just before the final (underlined) paragraph,
#    i.e. right after the last empty centered/bold paragraph that follows the chart.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">Time </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>comparison</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Comparing the original FP-Growth algorithm to our new proposed FP-Growth algorithm, we noticed some interesting performance differences. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">We compared these two algorithms using a dataset of 1, 2, 7, and 13 unique items on a database size of 10000, 50000, </w:t>
      </w:r>
      <w:r>
        <w:t>100000</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:t>500000</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, and </w:t>
      </w:r>
      <w:r>
        <w:t>1000000</w:t>
      </w:r>
      <w:r>
        <w:t>. We also experimented with a minimum support threshold of 3, and a much larger minimum support threshold of 10%</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Table </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>X</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">and table </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve">Y </w:t>
      </w:r>
      <w:r>
        <w:t>below illustrate the time (in seconds) that it took to execute each algorithm with a minimum support threshold equal to 10%</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and 3, respectively</w:t>
      </w:r>
      <w:r>
        <w:t>. Th</w:t>
      </w:r>
      <w:r>
        <w:t>ese</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> table</w:t>
      </w:r>
      <w:r>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>are</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> broken down into time it took to complete the first database scan, the time it took to complete the actual mining, and the total execution time for each algorithm.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>The time for the first database scan</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> in the original FP-Growth algorithm represents the time required to scan each transaction in the database and construct a simple list of transactions.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> With respect to our new FP-Growth algorithm, this time represents the time taken to scan each transaction in the database and construct our REP-Tree (a tree representing each transaction in the database</w:t>
      </w:r>
      <w:r>
        <w:t>).</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> Furthermore, the second column of data represents the time taken to complete the second database scan, determine the frequent 1-itemsets, and complete the mining process</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> to determine the remaining k-itemsets</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> with respect to the original FP-Growth algorithm.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> This column represents the time taken to iterate over the REP-Tree, determine the frequent 1-itemsets, and complete the mining process </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">to determine the remaining k-itemsets </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t>for our new proposed FP-Growth algorithm.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Looking at the results for each of these algorithms, we can see that the original FP-Growth algorithm </w:t>
      </w:r>
      <w:r>
        <w:t>outperforms</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> our new FP-Growth algorithm with large minimum support thresholds (see table </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>X</w:t>
      </w:r>
      <w:r>
        <w:t>).</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> This is largely due to the fact that the overhead to construct and iterate over the REP-Tree when a lot of items will inevitably be pruned due to not meeting the minimum support threshold requirement.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">It should be noted that our new FP-Growth algorithm performs all mining in a relatively constant time due to the tree traversal process. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">However, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">we see the opposite result when using a low minimum support threshold (see table </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Y</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">). </w:t>
      </w:r>
      <w:r>
        <w:t>In this scenario, the original FP-Growth algorithm must now iterate over almost all the possible itemsets, reading them off disk to construct the FP-Tree</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> While FP-Growth still </w:t>
      </w:r>
      <w:r>
        <w:t>outperforms</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> our new FP-Growth algorithm with on datasets with few unique items, </w:t>
      </w:r>
      <w:r>
        <w:t>our proposed algorithm</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> significantly outperforms the original FP-Growth algorithm when we begin mining on large datasets with many unique items.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>In terms of the practicality of our proposed changes to FP-Growth, these results are promising</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>Since</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> most interesting mining will be done using large datasets with many unique itemsets</w:t>
      </w:r>
      <w:r>
        <w:t>, our new algorithm could pose for more efficient frequent pattern mining</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> in certain scenarios.</w:t>
      </w:r>
    </w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

[void]$r.InsertXML($xml)
